$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values for the new timesheet rows. Shared-string entries are
# created in first-write order, so write the Task (column C) text in the
# same order the target shared-string table expects them:
#   "Data set corrections", "RapidFuzz parallel programming",
#   "RapidFuzz parallel programming and final matches data set", "Documentation"
$ws.Range("C46").Value = "Data set corrections"
$ws.Range("C43").Value = "RapidFuzz parallel programming"
$ws.Range("C44").Value = "RapidFuzz parallel programming"
$ws.Range("C45").Value = "RapidFuzz parallel programming and final matches data set"
$ws.Range("C47").Value = "Documentation"

$ws.Range("A43").Value = 44093
$ws.Range("B43").Value = 1

$ws.Range("A44").Value = 44094
$ws.Range("B44").Value = 1

$ws.Range("A45").Value = 44095
$ws.Range("B45").Value = 2

$ws.Range("A46").Value = 44097
$ws.Range("B46").Value = 1

$ws.Range("A47").Value = 44098
$ws.Range("B47").Value = 1

# Match the formatting of the preceding rows: columns A and B keep the same
# per-row style as row 42 (date style / numeric style), while column C picks
# up the plain numeric-style font (same as column B) instead of the task
# column's usual style.
$ws.Range("A42:B42").Copy()
$ws.Range("A43:B47").PasteSpecial(-4122)

$ws.Range("B42").Copy()
$ws.Range("C43:C47").PasteSpecial(-4122)

$ws.Range("C47").Select()
